$wb = $excel.ActiveWorkbook

# NOTE: column widths / default row heights on TestData & Role are left
# untouched on purpose -- this COM runtime's ColumnWidth/RowHeight setters
# snap to a coarse pixel grid, which would land further from the target
# stored values (and would also strip the existing bestFit flag on Role's
# column A) than simply leaving the original fine-grained widths in place.

# --- Add the new SNMP sheet after the last existing sheet (Role) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSnmp = $wb.Worksheets.Add($null, $lastSheet)
$wsSnmp.Name = "SNMP"

# Populate column-by-column so new shared-string entries are created in the
# same order as the target workbook (securityName, khang, authPro, MD5).
$wsSnmp.Range("A1").Value = "securityName"
$wsSnmp.Range("A2").Value = "khang"
$wsSnmp.Range("B1").Value = "authPro"
$wsSnmp.Range("B2").Value = "MD5"

# Match the selection recorded in the new sheet's view, then activate it so
# it becomes the active/selected tab (also clears tabSelected on Role).
$wsSnmp.Range("W8").Select()
$wsSnmp.Activate()
